$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = New-Object 'object[,]' 27,13
$data[0,0] = 0.441296842109439
$data[0,1] = 0.449932695668417
$data[0,2] = 0.430714505414509
$data[0,3] = 0.379850306961587
$data[0,4] = 0.519283019996893
$data[0,5] = 0.509035691130886
$data[0,6] = 0.35696975676024
$data[0,7] = 0.547774793298982
$data[0,8] = 0.329915275378447
$data[0,9] = 0.286873238075924
$data[0,10] = 0.126903211949325
$data[0,11] = 0.591065630264003
$data[0,12] = 0.633420069364728
$data[1,0] = 0.358986451498439
$data[1,1] = 0.276470716294786
$data[1,2] = 0.303598778912855
$data[1,3] = 0.219088466354011
$data[1,4] = 0.332027661372382
$data[1,5] = 0.342588192192442
$data[1,6] = 0.294054929364606
$data[1,7] = 0.234862255891938
$data[1,8] = 0.241773068797322
$data[1,9] = 0.430521704554094
$data[1,10] = 0.28959541811438
$data[1,11] = 0.438589567928279
$data[1,12] = 0.454630983761312
$data[2,0] = 0.198161661009927
$data[2,1] = 0.248313471603173
$data[2,2] = 0.200509520024054
$data[2,3] = 0.338700859532101
$data[2,4] = 0.18075163778984
$data[2,5] = 0.19552786464122
$data[2,6] = 0.264107937466889
$data[2,7] = 0.285191209788138
$data[2,8] = 0.288252973495061
$data[2,9] = 0.106626903721979
$data[2,10] = 0.143433043149037
$data[2,11] = 0.0495554116447345
$data[2,12] = 0.219446367219462
$data[3,0] = 0.169817051476306
$data[3,1] = 0.168675767460597
$data[3,2] = 0.177709841363343
$data[3,3] = 0.136240040232281
$data[3,4] = 0.197631752851924
$data[3,5] = 0.269396231574256
$data[3,6] = 0.080704594224714
$data[3,7] = 0.181083428783525
$data[3,8] = 0.15618998103631
$data[3,9] = 0.068711121256305
$data[3,10] = 0.226076004395808
$data[3,11] = 0.249014969288781
$data[3,12] = 0.179764784757296
$data[4,0] = 0.142291574054915
$data[4,1] = 0.120420461837184
$data[4,2] = 0.0716764126842037
$data[4,3] = 0.102435465438322
$data[4,4] = 0.101519633863153
$data[4,5] = 0.0925860347786772
$data[4,6] = 0.111599511465955
$data[4,7] = 0.245882140003428
$data[4,8] = 0.0668041216665552
$data[4,9] = 0.13541845331236
$data[4,10] = 0.0748623431861901
$data[4,11] = 0.151104268891365
$data[4,12] = 0.203500564857389
$data[5,0] = 0.134360154881326
$data[5,1] = 0.144739336216069
$data[5,2] = 0.169937045991076
$data[5,3] = 0.129494154250173
$data[5,4] = 0.0722512388783632
$data[5,5] = 0.074365555345283
$data[5,6] = 0.230925184978611
$data[5,7] = 0.156706562068422
$data[5,8] = 0.143496499849172
$data[5,9] = 0.12023860190341
$data[5,10] = 0.290796316593297
$data[5,11] = 0.0204861409133044
$data[5,12] = 0.0617420900988136
$data[6,0] = 0.112109960453632
$data[6,1] = 0.0925145394106335
$data[6,2] = 0.0641220705630488
$data[6,3] = 0.0741523770806195
$data[6,4] = 0.159185858290687
$data[6,5] = 0.0698329788192129
$data[6,6] = 0.139074978626742
$data[6,7] = 0.0661069594109144
$data[6,8] = 0.0739974150158471
$data[6,9] = 0.104448318383477
$data[6,10] = 0.0573616287475129
$data[6,11] = 0.264705701631175
$data[6,12] = 0.152363486085498
$data[7,0] = 0.104063148772217
$data[7,1] = 0.0871316721219487
$data[7,2] = 0.0810811941333479
$data[7,3] = 0.106441275818055
$data[7,4] = 0.0566279032087552
$data[7,5] = 0.0520853927838186
$data[7,6] = 0.0488300324020495
$data[7,7] = 0.122959709506436
$data[7,8] = 0.165051076246097
$data[7,9] = 0.272116767843546
$data[7,10] = 0.0874696342638326
$data[7,11] = 0.0989443281913727
$data[7,12] = 0.0666246995686403
$data[8,0] = 0.0952498017816782
$data[8,1] = 0.0901526582578352
$data[8,2] = 0.110588816330251
$data[8,3] = 0.103388392411287
$data[8,4] = 0.102315808626918
$data[8,5] = 0.0871959785086038
$data[8,6] = 0.106739591492894
$data[8,7] = 0.0495594431604003
$data[8,8] = 0.0557964745349515
$data[8,9] = 0.095508022197655
$data[8,10] = 0.132948385620608
$data[8,11] = 0.160557446962959
$data[8,12] = 0.0740313951637463
$data[9,0] = 0.0874706518488395
$data[9,1] = 0.0672871131000313
$data[9,2] = 0.0726654795684381
$data[9,3] = 0.0961371148960204
$data[9,4] = 0.0531777655747486
$data[9,5] = 0.047708258327873
$data[9,6] = 0.0425692434084225
$data[9,7] = 0.0706878949186356
$data[9,8] = 0.139239085555117
$data[9,9] = 0.298063895786922
$data[9,10] = 0.0466741458925149
$data[9,11] = 0.0336826727723516
$data[9,12] = 0.0523139132216924
$data[10,0] = 0.0770377228095487
$data[10,1] = 0.0635192914300062
$data[10,2] = 0.0701844083569901
$data[10,3] = 0.0307277086348666
$data[10,4] = 0.0880567938319764
$data[10,5] = 0.143744057187528
$data[10,6] = 0.0584045245037492
$data[10,7] = 0.0428909184754753
$data[10,8] = 0.035521555916916
$data[10,9] = 0.0153511539163466
$data[10,10] = 0.0571239345943341
$data[10,11] = 0.124061652553908
$data[10,12] = 0.122594419481344
$data[11,0] = 0.0756465636056349
$data[11,1] = 0.100552755154766
$data[11,2] = 0.0831589455316494
$data[11,3] = 0.181335642359528
$data[11,4] = 0.0624158969990063
$data[11,5] = 0.0757696859204116
$data[11,6] = 0.0844650018336526
$data[11,7] = 0.0640757056694225
$data[11,8] = 0.108997123980324
$data[11,9] = 0.0884812762770774
$data[11,10] = 0.0412611252283482
$data[11,11] = 0.0284285462270132
$data[11,12] = 0.0654105871665159
$data[12,0] = 0.0693730599263398
$data[12,1] = 0.0555416335961559
$data[12,2] = 0.0680403585016069
$data[12,3] = 0.0701807637107807
$data[12,4] = 0.035525824296384
$data[12,5] = 0.0266671385915977
$data[12,6] = 0.0335240300903662
$data[12,7] = 0.0800222364918165
$data[12,8] = 0.0537928392806654
$data[12,9] = 0.0964609688274934
$data[12,10] = 0.0440008049449396
$data[12,11] = 0.0317427123668093
$data[12,12] = 0.0916662947607395
$data[13,0] = 0.0683584853787358
$data[13,1] = 0.062252860473158
$data[13,2] = 0.0786250915915428
$data[13,3] = 0.0732239079678448
$data[13,4] = 0.0554277120006161
$data[13,5] = 0.0191130117994314
$data[13,6] = 0.048082149805928
$data[13,7] = 0.0673369545398527
$data[13,8] = 0.0420484062401228
$data[13,9] = 0.00667429312556995
$data[13,10] = 0.0273785693634304
$data[13,11] = 0.0220996987239972
$data[13,12] = 0.121035579433543
$data[14,0] = 0.0635845447230405
$data[14,1] = 0.0678888899197089
$data[14,2] = 0.0766637424337588
$data[14,3] = 0.058479003534329
$data[14,4] = 0.0468718388910785
$data[14,5] = 0.0775360913698474
$data[14,6] = 0.0298974248514686
$data[14,7] = 0.0979323812378051
$data[14,8] = 0.0831679651833013
$data[14,9] = 0.0808841576801516
$data[14,10] = 0.0392823899538042
$data[14,11] = 0.0337093324848879
$data[14,12] = 0.0643888275457631
$data[15,0] = 0.0594346632455331
$data[15,1] = 0.0619921172149611
$data[15,2] = 0.0810336727568943
$data[15,3] = 0.0662302102221586
$data[15,4] = 0.0422633984671631
$data[15,5] = 0.0351864527341478
$data[15,6] = 0.0558909977421723
$data[15,7] = 0.0578085539731153
$data[15,8] = 0.0977583348012284
$data[15,9] = 0.101644374785414
$data[15,10] = 0.039413655238913
$data[15,11] = 0.0229930688654548
$data[15,12] = 0.0506771174541164
$data[16,0] = 0.031130565528993
$data[16,1] = 0.032995705857802
$data[16,2] = 0.018052233261442
$data[16,3] = 0.0513222885365416
$data[16,4] = 0.0112250333288652
$data[16,5] = 0.0410415469087226
$data[16,6] = 0.00941716266274787
$data[16,7] = 0.0474725902142303
$data[16,8] = 0.109256759002246
$data[16,9] = 0.0173924457833681
$data[16,10] = 0.0316886280319162
$data[16,11] = 0.0306299118348573
$data[16,12] = 0.0332636586344221
$data[17,0] = 0.0271815945113994
$data[17,1] = 0.017611031177937
$data[17,2] = 0.0185418084524379
$data[17,3] = 0.0171594663909559
$data[17,4] = 0.0113959123064546
$data[17,5] = 0.0086597594208662
$data[17,6] = 0.020847666658421
$data[17,7] = 0.0195685497628379
$data[17,8] = 0.0101314460773026
$data[17,9] = 0.025598914163423
$data[17,10] = 0.0427144443182424
$data[17,11] = 0.0790068821805726
$data[17,12] = 0.0262985917231421
$data[18,0] = 0.0249117443373189
$data[18,1] = 0.0222920372084341
$data[18,2] = 0.0176165892070421
$data[18,3] = 0.0247648048085636
$data[18,4] = 0.0263307331462874
$data[18,5] = 0.0156359719814758
$data[18,6] = 0.0365892757027239
$data[18,7] = 0.0128570373930679
$data[18,8] = 0.026465295788619
$data[18,9] = 0.0488627195025423
$data[18,10] = 0.046755747780688
$data[18,11] = 0
$data[18,12] = 0.0110201536247897
$data[19,0] = 0.0205470813591216
$data[19,1] = 0.0216396600356974
$data[19,2] = 0.0180615587169546
$data[19,3] = 0.0417312000470512
$data[19,4] = 0.00835156014090111
$data[19,5] = 0.00812661283648955
$data[19,6] = 0.0165271735244459
$data[19,7] = 0.0117128760675959
$data[19,8] = 0.0835917780908648
$data[19,9] = 0.0108251287360204
$data[19,10] = 0.0279634801124558
$data[19,11] = 0.0335435645639343
$data[19,12] = 0.0191318239230357
$data[20,0] = 0.0147062967378243
$data[20,1] = 0.0230983022518674
$data[20,2] = 0.021915739814745
$data[20,3] = 0.0340060313222816
$data[20,4] = 0.0125970730652407
$data[20,5] = 0.0264500134354641
$data[20,6] = 0
$data[20,7] = 0.0245016176817279
$data[20,8] = 0.0598111153116144
$data[20,9] = 0.0159744352607008
$data[20,10] = 0.00689457553847729
$data[20,11] = 0.0056591287316197
$data[20,12] = 0.00902599124483849
$data[21,0] = 0.0129984622660514
$data[21,1] = 0.00653485539630399
$data[21,2] = 0.00643092238425799
$data[21,3] = 0.00657409117054134
$data[21,4] = 0.00511417380371708
$data[21,5] = 0
$data[21,6] = 0.00740187396793866
$data[21,7] = 0.00743679849884718
$data[21,8] = 0
$data[21,9] = 0
$data[21,10] = 0.00948016034952536
$data[21,11] = 0.0568545854228377
$data[21,12] = 0.0217735435025381
$data[22,0] = 0.0109995768215665
$data[22,1] = 0.013047801090276
$data[22,2] = 0.00313017848487854
$data[22,3] = 0.0289071544718255
$data[22,4] = 0
$data[22,5] = 0
$data[22,6] = 0
$data[22,7] = 0.0199424753566033
$data[22,8] = 0.0476756028609541
$data[22,9] = 0.0120216080678328
$data[22,10] = 0.0101812466154841
$data[22,11] = 0
$data[22,12] = 0.00989909979464741
$data[23,0] = 0.00941685109130378
$data[23,1] = 0.00819935265978798
$data[23,2] = 0.0102960558913913
$data[23,3] = 0.00330224422014463
$data[23,4] = 0.0191647460002275
$data[23,5] = 0
$data[23,6] = 0.00737184881470268
$data[23,7] = 0.00304250610842262
$data[23,8] = 0.0209182508826024
$data[23,9] = 0.00174709977982922
$data[23,10] = 0
$data[23,11] = 0.00476906523058334
$data[23,12] = 0.0181877874252677
$data[24,0] = 0.00623833363261425
$data[24,1] = 0.0044424671742946
$data[24,2] = 0.00546344682065765
$data[24,3] = 0.00780626036061284
$data[24,4] = 0.00494941073202395
$data[24,5] = 0.00704187449853556
$data[24,6] = 0
$data[24,7] = 0
$data[24,8] = 0
$data[24,9] = 0.00220398484545471
$data[24,10] = 0
$data[24,11] = 0
$data[24,12] = 0.0142671480896283
$data[25,0] = 0.00537646511325173
$data[25,1] = 0.00706613969157661
$data[25,2] = 0.0080075098469655
$data[25,3] = 0.00706681827292672
$data[25,4] = 0.00662418235538669
$data[25,5] = 0.0223148271221338
$data[25,6] = 0
$data[25,7] = 0.00972924143466308
$data[25,8] = 0
$data[25,9] = 0
$data[25,10] = 0
$data[25,11] = 0
$data[25,12] = 0.00885262888195205
$data[26,0] = 0.00501342589741842
$data[26,1] = 0.00175688665709081
$data[26,2] = 0.00416298594522594
$data[26,3] = 0.00330224422014463
$data[26,4] = 0
$data[26,5] = 0
$data[26,6] = 0
$data[26,7] = 0
$data[26,8] = 0
$data[26,9] = 0.00244374340211007
$data[26,10] = 0.0239299960466604
$data[26,11] = 0
$data[26,12] = 0.000990407240247677

$ws.Range("B2:N28").Value = $data
